$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H38").Value = 249.3
$ws.Range("I38").Value = 249.3
$ws.Range("K38").Value = 747.9000000000001
$ws.Range("M38").Value = -375.9000000000001

$ws.Range("H41").Value = 497.5484
$ws.Range("I41").Value = 547.93335
$ws.Range("J41").Value = 450.3125
$ws.Range("K41").Value = 547.93335
$ws.Range("L41").Value = 450.3125
$ws.Range("M41").Value = -107.93335
$ws.Range("N41").Value = -1330.3125

$ws.Range("H62").Value = 22246.215
$ws.Range("I62").Value = 18868.375
$ws.Range("K62").Value = 18868.375
$ws.Range("M62").Value = -18244.375

$ws.Range("H65").Value = 22246.215
$ws.Range("I65").Value = 18868.375
$ws.Range("K65").Value = 94341.875
$ws.Range("M65").Value = -91221.875

$ws.Range("H98").Value = 2210.9048
$ws.Range("J98").Value = 4548.6
$ws.Range("L98").Value = 4548.6
$ws.Range("N98").Value = -7544.6

$ws.Range("H106").Value = 40149.55
$ws.Range("I106").Value = 51687.75
$ws.Range("K106").Value = 51687.75
$ws.Range("M106").Value = -51056.75

$ws.Range("H122").Value = 2210.9048
$ws.Range("J122").Value = 4548.6
$ws.Range("L122").Value = 13645.8
$ws.Range("N122").Value = -18545.8

$ws.Range("H138").Value = 2831.1099
$ws.Range("J138").Value = 3263.3442
$ws.Range("L138").Value = 9790.0326
$ws.Range("N138").Value = -20070.0326

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H30").Value = 3003
$ws.Range("I30").Value = 9
$ws.Range("J30").Value = 4500
$ws.Range("K30").Value = 9
$ws.Range("L30").Value = 4500
$ws.Range("M30").Value = 141
$ws.Range("N30").Value = -4800

$ws.Range("H132").Value = 4146.787
$ws.Range("I132").Value = 4431.775
$ws.Range("J132").Value = 2518.2856
$ws.Range("K132").Value = 13295.325
$ws.Range("L132").Value = 7554.8568
$ws.Range("M132").Value = -10765.325
$ws.Range("N132").Value = -12614.8568

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H88").Value = 92177.60000000001
$ws.Range("J88").Value = 92177.60000000001
$ws.Range("L88").Value = 92177.60000000001
$ws.Range("N88").Value = -92989.60000000001

$ws.Range("H91").Value = 92177.60000000001
$ws.Range("J91").Value = 92177.60000000001
$ws.Range("L91").Value = 92177.60000000001
$ws.Range("N91").Value = -94985.60000000001

$ws.Range("H109").Value = 80000
$ws.Range("J109").Value = 80000
$ws.Range("L109").Value = 80000
$ws.Range("N109").Value = -82774

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 3985.6667
$ws.Range("J16").Value = 6256.5
$ws.Range("L16").Value = 6256.5
$ws.Range("N16").Value = -6830.5

$ws.Range("H22").Value = 198.43478
$ws.Range("I22").Value = 194.22223
$ws.Range("J22").Value = 213.6
$ws.Range("K22").Value = 194.22223
$ws.Range("L22").Value = 213.6
$ws.Range("M22").Value = 155.77777
$ws.Range("N22").Value = -913.6

$ws.Range("H31").Value = 3114.276
$ws.Range("I31").Value = 2269.611
$ws.Range("J31").Value = 4496.4546
$ws.Range("K31").Value = 2269.611
$ws.Range("L31").Value = 4496.4546
$ws.Range("M31").Value = -1974.611
$ws.Range("N31").Value = -5086.4546

$ws.Range("H34").Value = 3114.276
$ws.Range("I34").Value = 2269.611
$ws.Range("J34").Value = 4496.4546
$ws.Range("K34").Value = 2269.611
$ws.Range("L34").Value = 4496.4546
$ws.Range("M34").Value = -2067.611
$ws.Range("N34").Value = -4900.4546

$ws.Range("H39").Value = 1025.5
$ws.Range("I39").Value = 1025.5
$ws.Range("K39").Value = 1025.5
$ws.Range("M39").Value = -634.5

$ws.Range("H41").Value = 10999.333
$ws.Range("I41").Value = 5000
$ws.Range("J41").Value = 13999
$ws.Range("K41").Value = 5000
$ws.Range("L41").Value = 13999
$ws.Range("N41").Value = -14855
$ws.Range("M41").Value = -4572

$ws.Range("H49").Value = 1025.5
$ws.Range("I49").Value = 1025.5
$ws.Range("K49").Value = 1025.5
$ws.Range("M49").Value = -843.5

$ws.Range("H58").Value = 1522.75
$ws.Range("I58").Value = 1522.75
$ws.Range("K58").Value = 1522.75
$ws.Range("M58").Value = -1319.75

$ws.Range("H113").Value = 3985.6667
$ws.Range("J113").Value = 6256.5
$ws.Range("L113").Value = 6256.5
$ws.Range("N113").Value = -10596.5

$ws.Range("H136").Value = 1522.75
$ws.Range("I136").Value = 1522.75
$ws.Range("K136").Value = 4568.25
$ws.Range("M136").Value = -2018.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H37").Value = 125087180
$ws.Range("J37").Value = 125087180
$ws.Range("L37").Value = 375261540
$ws.Range("N37").Value = -375261764

$ws.Range("H56").Value = 10358.154
$ws.Range("I56").Value = 10358.154
$ws.Range("K56").Value = 10358.154
$ws.Range("M56").Value = -9828.154

$ws.Range("H93").Value = 336166.34
$ws.Range("J93").Value = 4250
$ws.Range("L93").Value = 12750
$ws.Range("N93").Value = -16494

$ws.Range("H129").Value = 2812.2273
$ws.Range("I129").Value = 1271.3334
$ws.Range("K129").Value = 3814.0002
$ws.Range("M129").Value = 1185.9998

$ws.Range("H134").Value = 2083.4546
$ws.Range("I134").Value = 2083.4546
$ws.Range("K134").Value = 6250.3638
$ws.Range("M134").Value = -1180.3638

$ws.Range("H139").Value = 70278.336
$ws.Range("I139").Value = 85625.914
$ws.Range("K139").Value = 256877.742
$ws.Range("M139").Value = -251737.742

$ws.Range("H140").Value = 2542.8
$ws.Range("I140").Value = 2329.1738
$ws.Range("J140").Value = 4999.5
$ws.Range("K140").Value = 6987.5214
$ws.Range("L140").Value = 14998.5
$ws.Range("M140").Value = -1807.5214
$ws.Range("N140").Value = -25358.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H92").Value = 21918.8
$ws.Range("J92").Value = 21918.8
$ws.Range("L92").Value = 21918.8
$ws.Range("N92").Value = -25662.8

$ws.Range("H136").Value = 62256.5
$ws.Range("J136").Value = 62256.5
$ws.Range("L136").Value = 186769.5
$ws.Range("N136").Value = -191869.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 3143
$ws.Range("I16").Value = 3143
$ws.Range("K16").Value = 3143
$ws.Range("M16").Value = -2973

$ws.Range("H22").Value = 1040.16
$ws.Range("I22").Value = 944.8570999999999
$ws.Range("K22").Value = 944.8570999999999
$ws.Range("M22").Value = -649.8570999999999

$ws.Range("H26").Value = 10500
$ws.Range("J26").Value = 10500
$ws.Range("L26").Value = 10500
$ws.Range("N26").Value = -11090

$ws.Range("H27").Value = 1040.16
$ws.Range("I27").Value = 944.8570999999999
$ws.Range("K27").Value = 944.8570999999999
$ws.Range("M27").Value = -837.8570999999999

$ws.Range("H58").Value = 6361
$ws.Range("I58").Value = 1736.6
$ws.Range("J58").Value = 14068.333
$ws.Range("K58").Value = 1736.6
$ws.Range("L58").Value = 14068.333
$ws.Range("M58").Value = -1476.6
$ws.Range("N58").Value = -14588.333

$ws.Range("H100").Value = 57944.176
$ws.Range("I100").Value = 24150.846
$ws.Range("K100").Value = 24150.846
$ws.Range("M100").Value = -23609.846

$ws.Range("H109").Value = 0
$ws.Range("I109").Value = 0
$ws.Range("J109").Value = 0
$ws.Range("K109").Value = 0
$ws.Range("L109").Value = 0
$ws.Range("M109").Value = ""
$ws.Range("N109").Value = ""

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H54").Value = 34582.168
$ws.Range("J54").Value = 46164.332
$ws.Range("L54").Value = 46164.332
$ws.Range("N54").Value = -47204.332

$ws.Range("H59").Value = 0
$ws.Range("J59").Value = 0
$ws.Range("L59").Value = 0
$ws.Range("N59").Value = ""

$ws.Range("H100").Value = 2724.0557
$ws.Range("I100").Value = 3103.3845
$ws.Range("J100").Value = 1737.8
$ws.Range("K100").Value = 6206.769
$ws.Range("L100").Value = 3475.6
$ws.Range("M100").Value = -5665.769
$ws.Range("N100").Value = -4557.6

$ws.Range("H132").Value = 5368.07
$ws.Range("I132").Value = 5732.4473
$ws.Range("K132").Value = 17197.3419
$ws.Range("M132").Value = -14667.3419
